$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "Shop Favs"
$ws.Range("F2").Value = "2025-03-30 18:20:06"

$ws.Range("E3").Value = "Amazon.com.br"
$ws.Range("F3").Value = "2025-03-30 18:20:07"

$ws.Range("E4").Value = "Mercado Livre Eletronicos"
$ws.Range("F4").Value = "2025-03-30 18:20:07"

$ws.Range("E5").Value = "SOMA3046530"
$ws.Range("F5").Value = "2025-03-30 18:20:08"
